# 20170601 Route 1 Fixes.
# Each departure time block originally repeated its row 5 times; trim each
# block down to 4 rows by deleting the 5th (last) occurrence in every
# block. Deleting from the bottom of the sheet upward keeps the remaining
# row numbers stable while we iterate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(93, 88, 83, 78, 73, 68, 63, 58, 53, 48, 43, 38, 33, 28, 23, 18, 13, 8)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
